# Update "Defect Report" sheet: fill in Functionality/Usability, Severity,
# Priority and Attachment (hyperlink) columns for rows 17-20, and move the
# active selection/viewport down to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect Report")

# xlPasteFormats = -4122 ; used to copy just the cell style (number format,
# font, fill, border, alignment) from an existing cell that already carries
# the style we need, without disturbing its value.
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 17 (Defect ID 18)
# ---------------------------------------------------------------------
$ws.Range("F17").Value = "Functionality,Usability"

$ws.Range("G2").Copy()
$ws.Range("G17").PasteSpecial($xlPasteFormats)
$ws.Range("G17").Value = "Medium"

$ws.Range("H2").Copy()
$ws.Range("H17").PasteSpecial($xlPasteFormats)
$ws.Range("H17").Value = "High"

$url17 = "https://drive.google.com/file/d/1pEoLGI23P4U2gYL0PslUmoLH1BAeP_Q1/view?usp=share_link"
$ws.Hyperlinks.Add($ws.Range("L17"), $url17, "", "", $url17)
$ws.Range("L2").Copy()
$ws.Range("L17").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(17).RowHeight = 105

# ---------------------------------------------------------------------
# Row 18 (Defect ID 28)
# ---------------------------------------------------------------------
$ws.Range("F18").Value = "Functionality,Usability"

$ws.Range("G3").Copy()
$ws.Range("G18").PasteSpecial($xlPasteFormats)
$ws.Range("G18").Value = "Low"

$ws.Range("H2").Copy()
$ws.Range("H18").PasteSpecial($xlPasteFormats)
$ws.Range("H18").Value = "Medium"

$url18 = "https://drive.google.com/file/d/1gQVo5aKPoCBtD9Uy23_07C8kcpZUu-aO/view?usp=share_link"
$ws.Hyperlinks.Add($ws.Range("L18"), $url18, "", "", $url18)
$ws.Range("L2").Copy()
$ws.Range("L18").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(18).RowHeight = 120

# ---------------------------------------------------------------------
# Row 19 (Defect ID 30)
# ---------------------------------------------------------------------
$ws.Range("F19").Value = "Functionality,Usability"

$ws.Range("G3").Copy()
$ws.Range("G19").PasteSpecial($xlPasteFormats)
$ws.Range("G19").Value = "Low"

$ws.Range("H3").Copy()
$ws.Range("H19").PasteSpecial($xlPasteFormats)
$ws.Range("H19").Value = "Low"

$url19 = "https://drive.google.com/file/d/1DlNEPMKd6O0Gj453c9VQnNaehj7pHxOR/view?usp=share_link"
$ws.Hyperlinks.Add($ws.Range("L19"), $url19, "", "", $url19)
$ws.Range("L2").Copy()
$ws.Range("L19").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(19).RowHeight = 120

# ---------------------------------------------------------------------
# Row 20 (Defect ID 33) - G20/H20 keep their original style (24)
# ---------------------------------------------------------------------
$ws.Range("F20").Value = "Functionality,Usability"
$ws.Range("G20").Value = "Low"
$ws.Range("H20").Value = "Low"

$url20 = "https://drive.google.com/file/d/1j-VDGXnzhtUdogdJxWvctwo9wZAd0Dz8/view?usp=share_link"
$ws.Hyperlinks.Add($ws.Range("L20"), $url20, "", "", $url20)
$ws.Range("L2").Copy()
$ws.Range("L20").PasteSpecial($xlPasteFormats)

$ws.Rows.Item(20).RowHeight = 120

# ---------------------------------------------------------------------
# Move viewport / selection to the new bottom row, as in the authored file
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A20").Select()

$excel.CutCopyMode = $false
